# feat: remove unused column from clo template
#
# The "expectedScorePercentage" column (originally column E) is unused and
# gets removed from the CLO template: the worksheet column is deleted (so
# everything to its right shifts left), the backing table is resized down
# to match, and the table's column names are re-synced from the (now
# shifted) header row text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the column to drop by its header text so this keeps working even
# if the sheet layout ever changes.
$headerRow = 1
$lastCol = $ws.UsedRange.Columns.Count
$targetColIndex = -1
for ($c = 1; $c -le $lastCol; $c++) {
    $headerValue = $ws.Cells.Item($headerRow, $c).Value()
    if ($headerValue -eq "expectedScorePercentage") {
        $targetColIndex = $c
    }
}

if ($targetColIndex -eq -1) {
    throw "Could not find the 'expectedScorePercentage' column to remove"
}

# Remember the table so we can resize/resync it after the column shift.
$tbl = $ws.ListObjects.Item(1)
$origTopLeft = $tbl.Range.Cells.Item(1, 1)
$origRowCount = $tbl.Range.Rows.Count
$newColCount = $tbl.ListColumns.Count - 1

# Delete the whole worksheet column - shifts every column to its right one
# position to the left, fixing up cell references/shared strings as it goes.
$ws.Columns.Item($targetColIndex).Delete()

# Shrink the table definition to match the new (narrower) data range.
$newRange = $ws.Range($origTopLeft, $ws.Cells.Item($origTopLeft.Row + $origRowCount - 1, $origTopLeft.Column + $newColCount - 1))
$tbl.Resize($newRange)

# Re-assign each header cell's own value so the table's column names
# re-sync with the (now shifted) header text instead of staying stale.
for ($c = 1; $c -le $newColCount; $c++) {
    $headerCell = $ws.Cells.Item($headerRow, $c)
    $headerCell.Value = $headerCell.Value()
}

# Match the author's final cursor position.
$ws.Range("E11").Select() | Out-Null
